{"js": "// Replace the date line and each \"A\u00d7B=C\" multiplication-table answer cell\n// with its updated value (26 total text replacements: 1 date + 25 table cells).\nconst replacements = [\n  [\"2024-10-30 Wednesday\", \"2024-10-31 Thursday\"],\n  [\"844\u00d73=2532\", \"174\u00d74=696\"],\n  [\"936\u00d75=4680\", \"747\u00d72=1494\"],\n  [\"545\u00d79=4905\", \"659\u00d74=2636\"],\n  [\"412\u00d72=824\", \"693\u00d74=2772\"],\n  [\"584\u00d77=4088\", \"745\u00d77=5215\"],\n  [\"875\u00d74=3500\", \"139\u00d72=278\"],\n  [\"211\u00d77=1477\", \"951\u00d72=1902\"],\n  [\"506\u00d76=3036\", \"145\u00d74=580\"],\n  [\"262\u00d77=1834\", \"905\u00d79=8145\"],\n  [\"537\u00d78=4296\", \"963\u00d74=3852\"],\n  [\"471\u00d74=1884\", \"600\u00d76=3600\"],\n  [\"994\u00d76=5964\", \"561\u00d76=3366\"],\n  [\"746\u00d78=5968\", \"428\u00d72=856\"],\n  [\"579\u00d79=5211\", \"511\u00d74=2044\"],\n  [\"640\u00d72=1280\", \"329\u00d72=658\"],\n  [\"837\u00d72=1674\", \"813\u00d77=5691\"],\n  [\"407\u00d77=2849\", \"682\u00d74=2728\"],\n  [\"620\u00d74=2480\", \"611\u00d72=1222\"],\n  [\"853\u00d75=4265\", \"568\u00d74=2272\"],\n  [\"157\u00d74=628\", \"981\u00d75=4905\"],\n  [\"690\u00d73=2070\", \"171\u00d76=1026\"],\n  [\"662\u00d73=1986\", \"333\u00d76=1998\"],\n  [\"921\u00d72=1842\", \"743\u00d74=2972\"],\n  [\"608\u00d72=1216\", \"449\u00d77=3143\"],\n  [\"132\u00d73=396\", \"758\u00d76=4548\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00d7B=C\" multiplication-table answer cell\n# with its updated value (26 total text replacements: 1 date + 25 table cells).\n$d = $word.ActiveDocument\n\n$oldValues = @(\n  \"2024-10-30 Wednesday\",\n  \"844\u00d73=2532\",\n  \"936\u00d75=4680\",\n  \"545\u00d79=4905\",\n  \"412\u00d72=824\",\n  \"584\u00d77=4088\",\n  \"875\u00d74=3500\",\n  \"211\u00d77=1477\",\n  \"506\u00d76=3036\",\n  \"262\u00d77=1834\",\n  \"537\u00d78=4296\",\n  \"471\u00d74=1884\",\n  \"994\u00d76=5964\",\n  \"746\u00d78=5968\",\n  \"579\u00d79=5211\",\n  \"640\u00d72=1280\",\n  \"837\u00d72=1674\",\n  \"407\u00d77=2849\",\n  \"620\u00d74=2480\",\n  \"853\u00d75=4265\",\n  \"157\u00d74=628\",\n  \"690\u00d73=2070\",\n  \"662\u00d73=1986\",\n  \"921\u00d72=1842\",\n  \"608\u00d72=1216\",\n  \"132\u00d73=396\"\n)\n$newValues = @(\n  \"2024-10-31 Thursday\",\n  \"174\u00d74=696\",\n  \"747\u00d72=1494\",\n  \"659\u00d74=2636\",\n  \"693\u00d74=2772\",\n  \"745\u00d77=5215\",\n  \"139\u00d72=278\",\n  \"951\u00d72=1902\",\n  \"145\u00d74=580\",\n  \"905\u00d79=8145\",\n  \"963\u00d74=3852\",\n  \"600\u00d76=3600\",\n  \"561\u00d76=3366\",\n  \"428\u00d72=856\",\n  \"511\u00d74=2044\",\n  \"329\u00d72=658\",\n  \"813\u00d77=5691\",\n  \"682\u00d74=2728\",\n  \"611\u00d72=1222\",\n  \"568\u00d74=2272\",\n  \"981\u00d75=4905\",\n  \"171\u00d76=1026\",\n  \"333\u00d76=1998\",\n  \"743\u00d74=2972\",\n  \"449\u00d77=3143\",\n  \"758\u00d76=4548\"\n)\n\nfor ($i = 0; $i -lt $oldValues.Count; $i++) {\n  $oldText = $oldValues[$i]\n  $newText = $newValues[$i]\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n"}
